# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Membrillo" sheet.
# The data rows (2-35) are re-ordered (re-sorted by date) while the set of
# records itself is unchanged, so we snapshot every data row first and then
# write each source row's values into its new destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 35
$firstCol = 1
$lastCol = 20

# Mapping of destination row -> source row (1-indexed worksheet rows).
$rowMap = @{
    2  = 28
    3  = 24
    4  = 16
    5  = 17
    6  = 19
    7  = 20
    8  = 21
    9  = 4
    10 = 5
    11 = 22
    12 = 23
    13 = 10
    14 = 11
    15 = 12
    16 = 32
    17 = 33
    18 = 2
    19 = 3
    20 = 13
    21 = 14
    22 = 15
    23 = 6
    24 = 27
    25 = 7
    26 = 8
    27 = 9
    28 = 29
    29 = 30
    30 = 31
    31 = 26
    32 = 25
    33 = 34
    34 = 35
    35 = 18
}

# Snapshot every source row (full A:T range) before any writes happen, since
# a row can be simultaneously a source for one destination and a destination
# for another.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# Write the snapshotted rows back out in their new order.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $data = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $data[$c - $firstCol]
    }
}
